$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 77.30653456323489
$ws.Range("B3").Value = 0.8937246874339718
$ws.Range("B4").Value = 0.06148432582916463
$ws.Range("B5").Value = 0.3920582124178325
